$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.985.53"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.18"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.90%  "
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4308"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3697"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07252"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.160.42"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +22.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8701"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.33"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.95%  "
$ws.Range("E13").Value = "  +4.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.420"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06973"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.11"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.015"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008927"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.023.55"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.214"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.41%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.99"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.376.02"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +19.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.77"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.889"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.37"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.246"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.930"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +14.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.96"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08977"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("E32").Value = "  +6.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7458"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.440"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.818"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.009"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.125"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05241"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01929"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5122"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.754"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +10.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1654"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.519"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.313"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "107.57"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.43"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.009"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.652"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.36%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06308"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4572"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.816"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.96%  "
